$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44371
$ws.Range("J2").Value = 150
$ws.Range("D3").Value = 44364
$ws.Range("J3").Value = 100
$ws.Range("O3").Value = 'Región Metropolitana'
$ws.Range("D4").Value = 44369
$ws.Range("N4").Value = '$/caja 20 docenas'
$ws.Range("P4").Value = 7000
$ws.Range("Q4").Value = 1
$ws.Range("D5").Value = 44354
$ws.Range("O5").Value = 'Región del Maule'
$ws.Range("D6").Value = 44355
$ws.Range("O6").Value = 'Región Metropolitana'
$ws.Range("D7").Value = 44358
$ws.Range("J7").Value = 150
$ws.Range("K7").Value = 7000
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 7000
$ws.Range("P7").Value = 194
$ws.Range("D8").Value = 44386
$ws.Range("J8").Value = 200
$ws.Range("K8").Value = 6500
$ws.Range("L8").Value = 6500
$ws.Range("M8").Value = 6500
$ws.Range("P8").Value = 181
$ws.Range("D9").Value = 44357
$ws.Range("K9").Value = 6500
$ws.Range("L9").Value = 6500
$ws.Range("M9").Value = 6500
$ws.Range("N9").Value = '$/caja 20 docenas'
$ws.Range("P9").Value = 6500
$ws.Range("Q9").Value = 1
$ws.Range("D10").Value = 44342
$ws.Range("J10").Value = 150
$ws.Range("N10").Value = '$/caja 36 atados'
$ws.Range("O10").Value = 'Región del Maule'
$ws.Range("P10").Value = 194
$ws.Range("Q10").Value = 36
$ws.Range("D11").Value = 44372
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("P11").Value = 194
$ws.Range("D12").Value = 44376
$ws.Range("N12").Value = '$/caja 36 atados'
$ws.Range("O12").Value = 'Región Metropolitana'
$ws.Range("P12").Value = 181
$ws.Range("Q12").Value = 36
$ws.Range("D13").Value = 44362
$ws.Range("J13").Value = 100
$ws.Range("D14").Value = 44348
$ws.Range("O14").Value = 'Región del Maule'
$ws.Range("D15").Value = 44340
